$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2193.8125
$ws.Range("I40").Value = 2533.3333
$ws.Range("J40").Value = 1990.1
$ws.Range("K40").Value = 2533.3333
$ws.Range("L40").Value = 1990.1
$ws.Range("M40").Value = -2358.3333
$ws.Range("N40").Value = -2340.1

# Row 116
$ws.Range("H116").Value = 3942.8572
$ws.Range("I116").Value = 3942.8572
$ws.Range("K116").Value = 3942.8572
$ws.Range("M116").Value = -500.8571999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 9472
$ws.Range("J43").Value = 9472
$ws.Range("L43").Value = 9472
$ws.Range("N43").Value = -10098

# Row 61
$ws.Range("H61").Value = 3278
$ws.Range("I61").Value = 2060.5
$ws.Range("K61").Value = 2060.5
$ws.Range("M61").Value = -1848.5

# Row 74
$ws.Range("H74").Value = 7475.6
$ws.Range("I74").Value = 1233.1666
$ws.Range("J74").Value = 16839.25
$ws.Range("K74").Value = 1233.1666
$ws.Range("L74").Value = 16839.25
$ws.Range("M74").Value = -359.1666
$ws.Range("N74").Value = -18587.25

# Row 77
$ws.Range("H77").Value = 7475.6
$ws.Range("I77").Value = 1233.1666
$ws.Range("J77").Value = 16839.25
$ws.Range("K77").Value = 6165.833000000001
$ws.Range("L77").Value = 84196.25
$ws.Range("M77").Value = -1797.833000000001
$ws.Range("N77").Value = -92932.25

# Row 122
$ws.Range("H122").Value = 4390.696
$ws.Range("I122").Value = 3922
$ws.Range("K122").Value = 11766
$ws.Range("M122").Value = -9316

# Row 136
$ws.Range("H136").Value = 3278
$ws.Range("I136").Value = 2060.5
$ws.Range("K136").Value = 6181.5
$ws.Range("M136").Value = -3631.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 496.90323
$ws.Range("I94").Value = 379.16
$ws.Range("J94").Value = 987.5
$ws.Range("K94").Value = 379.16
$ws.Range("L94").Value = 987.5
$ws.Range("M94").Value = 71.83999999999997
$ws.Range("N94").Value = -1889.5

# Row 105
$ws.Range("H105").Value = 246751.69
$ws.Range("I105").Value = 2636.2334
$ws.Range("K105").Value = 2636.2334
$ws.Range("M105").Value = -889.2334000000001

# Row 107
$ws.Range("H107").Value = 902.75
$ws.Range("I107").Value = 370.33334
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 370.33334
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 1549.66666
$ws.Range("N107").Value = -6340

# Row 134
$ws.Range("H134").Value = 3393.513
$ws.Range("I134").Value = 2354.8462
$ws.Range("J134").Value = 5470.846
$ws.Range("K134").Value = 7064.5386
$ws.Range("L134").Value = 16412.538
$ws.Range("M134").Value = -4529.5386
$ws.Range("N134").Value = -21482.538

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 41669572
$ws.Range("I86").Value = 83335544
$ws.Range("K86").Value = 83335544
$ws.Range("M86").Value = -83334421

# Row 89
$ws.Range("H89").Value = 41669572
$ws.Range("I89").Value = 83335544
$ws.Range("K89").Value = 416677720
$ws.Range("M89").Value = -416672104

# Row 99
$ws.Range("H99").Value = 4983.75
$ws.Range("I99").Value = 6202.8125
$ws.Range("J99").Value = 3358.3333
$ws.Range("K99").Value = 6202.8125
$ws.Range("L99").Value = 3358.3333
$ws.Range("M99").Value = -4704.8125
$ws.Range("N99").Value = -6354.3333

# Row 126
$ws.Range("H126").Value = 4983.75
$ws.Range("I126").Value = 6202.8125
$ws.Range("J126").Value = 3358.3333
$ws.Range("K126").Value = 18608.4375
$ws.Range("L126").Value = 10074.9999
$ws.Range("M126").Value = -16138.4375
$ws.Range("N126").Value = -15014.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1633.0416
$ws.Range("J34").Value = 1567.6818
$ws.Range("L34").Value = 4703.0454
$ws.Range("N34").Value = -4871.0454

# Row 39
$ws.Range("H39").Value = 8148.276
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 8148.276
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 24444.828
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -25032.828

# Row 55
$ws.Range("H55").Value = 3021.818
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 3224
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 9672
$ws.Range("N55").Value = -10026
$ws.Range("M55").Value = -2823

# Row 108
$ws.Range("H108").Value = 1782.625
$ws.Range("I108").Value = 1180.1428
$ws.Range("J108").Value = 6000
$ws.Range("K108").Value = 3540.4284
$ws.Range("L108").Value = 18000
$ws.Range("M108").Value = -660.4284000000002
$ws.Range("N108").Value = -23760

# Row 119
$ws.Range("H119").Value = 1061.4445
$ws.Range("I119").Value = 564.7143
$ws.Range("K119").Value = 1694.1429
$ws.Range("M119").Value = 3143.8571

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2476.182
$ws.Range("I122").Value = 2285.7144
$ws.Range("J122").Value = 2809.5
$ws.Range("K122").Value = 6857.1432
$ws.Range("L122").Value = 8428.5
$ws.Range("M122").Value = -4407.1432
$ws.Range("N122").Value = -13328.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 907.5
$ws.Range("I22").Value = 640
$ws.Range("J22").Value = 1308.75
$ws.Range("K22").Value = 640
$ws.Range("L22").Value = 1308.75
$ws.Range("M22").Value = -345
$ws.Range("N22").Value = -1898.75

# Row 27
$ws.Range("H27").Value = 907.5
$ws.Range("I27").Value = 640
$ws.Range("J27").Value = 1308.75
$ws.Range("K27").Value = 640
$ws.Range("L27").Value = 1308.75
$ws.Range("M27").Value = -533
$ws.Range("N27").Value = -1522.75

# Row 93
$ws.Range("H93").Value = 2419.3333
$ws.Range("I93").Value = 2459.8
$ws.Range("J93").Value = 2390.4285
$ws.Range("K93").Value = 2459.8
$ws.Range("L93").Value = 2390.4285
$ws.Range("M93").Value = -1211.8
$ws.Range("N93").Value = -4886.4285

# Row 136
$ws.Range("H136").Value = 4761.3335
$ws.Range("I136").Value = 2354.889
$ws.Range("J136").Value = 19200
$ws.Range("K136").Value = 7064.667
$ws.Range("L136").Value = 57600
$ws.Range("M136").Value = -4514.667
$ws.Range("N136").Value = -62700

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 6081.3335
$ws.Range("J41").Value = 6229.2
$ws.Range("L41").Value = 6229.2
$ws.Range("N41").Value = -7009.2

# Row 107
$ws.Range("H107").Value = 4126.5454
$ws.Range("I107").Value = 5349
$ws.Range("J107").Value = 866.6667
$ws.Range("K107").Value = 16047
$ws.Range("L107").Value = 2600.0001
$ws.Range("M107").Value = -14127
$ws.Range("N107").Value = -6440.0001

# Row 136
$ws.Range("H136").Value = 1521.44
$ws.Range("I136").Value = 848.82355
$ws.Range("J136").Value = 2950.75
$ws.Range("K136").Value = 2546.47065
$ws.Range("L136").Value = 8852.25
$ws.Range("M136").Value = 3.52935000000025
$ws.Range("N136").Value = -13952.25
